# #987 - runtime: rethink core struct
#
# Applies:
#   1) overflowPunct: false -> true on the title paragraph's pPr
#   2) overflowPunct: false -> true on the "fixnums" section paragraph's pPr
#   3) overflowPunct: false -> true on the Normal style's pPr (styles.xml)
#   4) Text-run fix: "[dependenc" + "(featu" + "ies]" -> "[featu" + "r" + "es]"
#      (so the visible text "[dependenc(featuies]" becomes "[features]")

$d = $word.ActiveDocument

# --- 1) & 2): overflowPunct false -> true -------------------------------
# ParagraphFormat.HangingPunctuation is the VBA property backing
# w:overflowPunct (Format > Asian Typography > "Allow punctuation marks to
# extend past the text margin"). Flip it for the two paragraphs whose pPr
# currently carries <w:overflowPunct w:val="false"/>.

$n = $d.Paragraphs.Count
for ($i = 1; $i -le $n; $i++) {
    $p = $d.Paragraphs($i)
    $t = $p.Range.Text
    if ($t.Contains("Mu Runtime Reference") -or $t.Contains("fixnums")) {
        $p.Format.HangingPunctuation = $true
    }
}

# --- 3) overflowPunct false -> true on the Normal style ------------------
$normal = $d.Styles("Normal")
$normal.ParagraphFormat.HangingPunctuation = $true

# --- 4) Fix the split "[dependenc(featuies]" run text --------------------
# Find the paragraph that holds the broken run text.
$target = $null
for ($i = 1; $i -le $n; $i++) {
    $p = $d.Paragraphs($i)
    if ($p.Range.Text.Contains("[dependenc(featuies]")) {
        $target = $p
        break
    }
}

if ($target -ne $null) {
    $pStart = $target.Range.Start

    # Known layout relative to the paragraph start:
    #   +0 .. +3  "   "          (left untouched)
    #   +3 .. +13 "[dependenc"   -> "[featu"
    #   +13 .. +19 "(featu"      -> "r"
    #   +19 .. +23 "ies]"        -> "es]"
    $s0 = $pStart
    $s1 = $pStart + 3
    $s2 = $pStart + 13
    $s3 = $pStart + 19
    $s4 = $pStart + 23

    # Replace right-to-left so untouched offsets to the left stay valid.
    $r4 = $d.Range($s3, $s4)
    $r4.Find.Execute("ies]", $false, $false, $false, $false, $false, $true, 1, $false, "es]", 2) | Out-Null

    $r3 = $d.Range($s2, $s3)
    $r3.Find.Execute("(featu", $false, $false, $false, $false, $false, $true, 1, $false, "r", 2) | Out-Null

    $r2 = $d.Range($s1, $s2)
    $r2.Find.Execute("[dependenc", $false, $false, $false, $false, $false, $true, 1, $false, "[featu", 2) | Out-Null

    # The engine coalesces same-format runs touched by an edit into a
    # single run. Re-split it back into the original run boundaries
    # (one run per original <w:r>) by toggling Bold off/on at each
    # boundary - a no-op on the visible formatting, but it forces the
    # writer to keep the runs distinct.
    $newEnd = $s1 + 6 + 1 + 3   # "[featu" + "r" + "es]" = 10 chars after s1

    $d.Range($s1, $newEnd).Font.Bold = $false
    $d.Range($s1, $newEnd).Font.Bold = $true

    $d.Range($s1 + 6, $newEnd).Font.Bold = $false
    $d.Range($s1 + 6, $newEnd).Font.Bold = $true

    $d.Range($s1 + 7, $newEnd).Font.Bold = $false
    $d.Range($s1 + 7, $newEnd).Font.Bold = $true
}

Write-Output "done"
